$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.609.42"
$ws.Range("E2").Value = "  -1.40%  "

# Row 3
$ws.Range("D3").Value = "3.067.36"
$ws.Range("E3").Value = "  -5.12%  "

# Row 4
$ws.Range("E4").Value = "  -0.52%  "

# Row 5
$ws.Range("D5").Value = "'586.02"
$ws.Range("E5").Value = "  -1.61%  "

# Row 6
$ws.Range("D6").Value = "'154.16"
$ws.Range("E6").Value = "  +3.40%  "

# Row 7
$ws.Range("E7").Value = "  -0.43%  "

# Row 8
$ws.Range("D8").Value = "'0.535"
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").Value = "3.065.37"
$ws.Range("E9").Value = "  -3.60%  "

# Row 10
$ws.Range("E10").Value = "  -4.96%  "

# Row 11
$ws.Range("D11").Value = "'5.81"
$ws.Range("E11").Value = "  -4.35%  "

# Row 12
$ws.Range("E12").Value = "  -3.61%  "

# Row 13
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "'0.0000236"
$ws.Range("E13").Value = "  -5.08%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'36.59"
$ws.Range("E14").Value = "  -3.72%  "

# Row 15
$ws.Range("E15").Value = "  -2.22%  "

# Row 16
$ws.Range("D16").Value = "3.575.44"
$ws.Range("E16").Value = "  -5.06%  "

# Row 17
$ws.Range("D17").Value = "63.536.28"
$ws.Range("E17").Value = "  -1.21%  "

# Row 18
$ws.Range("E18").Value = "  -3.60%  "

# Row 19
$ws.Range("D19").Value = "3.069.10"
$ws.Range("E19").Value = "  -3.99%  "

# Row 20
$ws.Range("D20").Value = "'469.53"
$ws.Range("E20").Value = "  -1.44%  "

# Row 21
$ws.Range("D21").Value = "'14.23"
$ws.Range("E21").Value = "  -2.98%  "

# Row 22
$ws.Range("E22").Value = "  -5.72%  "

# Row 23
$ws.Range("E23").Value = "  -3.37%  "

# Row 24
$ws.Range("D24").Value = "'2.44"
$ws.Range("E24").Value = "  -1.34%  "

# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'12.75"
$ws.Range("E25").Value = "  -4.60%  "

# Row 26
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'80.17"
$ws.Range("E26").Value = "  -2.21%  "

# Row 27
$ws.Range("D27").Value = "'10.41"
$ws.Range("E27").Value = "  +4.04%  "

# Row 28
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.23%  "

# Row 29
$ws.Range("E29").Value = "  +1.21%  "

# Row 30
$ws.Range("E30").Value = "  -0.41%  "

# Row 31
$ws.Range("E31").Value = "  -3.46%  "

# Row 32
$ws.Range("E32").Value = "  -6.56%  "

# Row 33
$ws.Range("E33").Value = "  -8.67%  "

# Row 34
$ws.Range("D34").Value = "'26.96"
$ws.Range("E34").Value = "  -5.27%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0820"
$ws.Range("E35").Value = "  -5.29%  "

# Row 36
$ws.Range("E36").Value = "  -3.14%  "

# Row 37
$ws.Range("D37").Value = "'5.95"
$ws.Range("E37").Value = "  -5.51%  "

# Row 38
$ws.Range("D38").Value = "'3.24"
$ws.Range("E38").Value = "  -3.93%  "

# Row 39
$ws.Range("D39").Value = "'2.20"
$ws.Range("E39").Value = "  -5.69%  "

# Row 40
$ws.Range("D40").Value = "'50.47"
$ws.Range("E40").Value = "  -2.62%  "

# Row 41
$ws.Range("D41").Value = "'9.09"
$ws.Range("E41").Value = "  -3.62%  "

# Row 42
$ws.Range("D42").Value = "'436.66"
$ws.Range("E42").Value = "  -7.20%  "

# Row 43
$ws.Range("E43").Value = "  -4.03%  "

# Row 44
$ws.Range("D44").Value = "'40.41"
$ws.Range("E44").Value = "  +2.23%  "

# Row 45
$ws.Range("E45").Value = "  +0.55%  "

# Row 46
$ws.Range("E46").Value = "  -5.41%  "

# Row 47
$ws.Range("D47").Value = "2.794.45"
$ws.Range("E47").Value = "  -5.18%  "

# Row 48
$ws.Range("D48").Value = "'130.22"
$ws.Range("E48").Value = "  -1.96%  "

# Row 49
$ws.Range("E49").Value = "  +0.08%  "

# Row 50
$ws.Range("D50").Value = "'24.85"
$ws.Range("E50").Value = "  +0.00%  "

# Row 51
$ws.Range("D51").Value = "'2.20"
$ws.Range("E51").Value = "  -3.86%  "
